$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.398.53"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").Value = "1.941.33"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.54"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.91"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -3.58%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.365"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "55.84"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0837"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +7.69%  "
$ws.Range("E12").Value = "  +1.21%  "
$ws.Range("E13").Value = "  -2.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.49"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.80%  "
$ws.Range("D15").Value = "2.226.04"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.60"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("E17").Value = "  -1.65%  "
$ws.Range("D18").Value = "1.943.00"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "36.336.16"
$ws.Range("E19").Value = "  +2.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.70"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("E21").Value = "  +3.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.44"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.05"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.49%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.43"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.27"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.10"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.40"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.124"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.47%  "
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("E32").Value = "  +3.03%  "
$ws.Range("E33").Value = "  -2.33%  "
$ws.Range("E34").Value = "  +3.44%  "
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.23"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.78"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.31%  "
$ws.Range("E39").Value = "  -5.45%  "
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.93"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.84%  "
$ws.Range("E43").Value = "  -2.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0209"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.08"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.95%  "
$ws.Range("D46").Value = "1.353.21"
$ws.Range("E46").Value = "  +2.23%  "
$ws.Range("E47").Value = "  -3.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.65"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.14"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.93%  "
$ws.Range("E50").Value = "  +3.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.26"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.91%  "
